# Update the evaluation form's title heading.
#
# Before: "Lakbay: Evaluation Form for Student Drivers/Citizens"
# After:  "Consolidated Sample Response of a Student Driver/Citizen"
#
# (The title paragraph keeps its Heading1 style; only the run text changes.)

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Lakbay: Evaluation Form for Student Drivers/Citizens",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Consolidated Sample Response of a Student Driver/Citizen",
    2
) | Out-Null
